# Fix the query text in B3: ORDER By samp.sample_id LIMIT 10 -> LIMIT 100
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("B3")
$current = $range.Value2
$updated = $current -replace 'ORDER By samp\.sample_id LIMIT 10$', 'ORDER By samp.sample_id LIMIT 100'
$range.Value = $updated

# Move / update the active selection to B13 (as captured in the saved view state)
$ws.Range("B13").Select()
